# Swap the taxon/observation data between row 13 and row 14, while
# leaving the shared/location/observer columns untouched in each row.
#
# Row 13 currently holds the "Kandelabersvamp / Artomyces pyxidatus" record
# and row 14 currently holds the "Apelticka / Aurantiporus fissilis" record.
# After the edit, row 13 should hold the Apelticka record (including its
# public-comment and substrate columns) and row 14 should hold the
# Kandelabersvamp record (with those extra columns cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with the taxon record between the
# two rows.
$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $v13 = $ws.Range("$col" + "13").Value()
    $v14 = $ws.Range("$col" + "14").Value()
    $ws.Range("$col" + "13").Value = $v14
    $ws.Range("$col" + "14").Value = $v13
}

# AC13/AJ13/AK13/AO13 were empty and now receive the comment/substrate
# values that used to live on row 14; the row-14 copies are cleared out.
$ws.Range("AC13").Value = "2 fruktkroppar på två träd ett stående dött träd och en låga."
$ws.Range("AJ13").Value = "asp"
$ws.Range("AK13").Value = "Populus tremula"
$ws.Range("AO13").Value = "Populus tremula"

$ws.Range("AC14").ClearContents()
$ws.Range("AJ14").ClearContents()
$ws.Range("AK14").ClearContents()
$ws.Range("AO14").ClearContents()
